$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 7396
$ws.Range("I38").Value = 7396
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 22188
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -21816
$ws.Range("N38").ClearContents()
$ws.Range("H64").Value = 4625
$ws.Range("J64").Value = 4833.3335
$ws.Range("L64").Value = 4833.3335
$ws.Range("N64").Value = -5329.3335
$ws.Range("H67").Value = 4625
$ws.Range("J67").Value = 4833.3335
$ws.Range("L67").Value = 4833.3335
$ws.Range("N67").Value = -6549.3335
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H100").Value = 1084.7142
$ws.Range("I100").Value = 1007.36365
$ws.Range("K100").Value = 1007.36365
$ws.Range("M100").Value = -466.36365
$ws.Range("H113").Value = 3689.8
$ws.Range("I113").Value = 2950
$ws.Range("J113").Value = 3874.75
$ws.Range("K113").Value = 2950
$ws.Range("L113").Value = 3874.75
$ws.Range("M113").Value = 304
$ws.Range("N113").Value = -10382.75
$ws.Range("H132").Value = 2551.3076
$ws.Range("I132").Value = 2670.6365
$ws.Range("K132").Value = 8011.9095
$ws.Range("M132").Value = -5481.9095
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2280.2942
$ws.Range("I61").Value = 2277.7334
$ws.Range("K61").Value = 2277.7334
$ws.Range("M61").Value = -2065.7334
$ws.Range("H74").Value = 1372.4445
$ws.Range("I74").Value = 1911.3334
$ws.Range("J74").Value = 294.66666
$ws.Range("K74").Value = 1911.3334
$ws.Range("L74").Value = 294.66666
$ws.Range("M74").Value = -1037.3334
$ws.Range("N74").Value = -2042.66666
$ws.Range("H77").Value = 1372.4445
$ws.Range("I77").Value = 1911.3334
$ws.Range("J77").Value = 294.66666
$ws.Range("K77").Value = 9556.666999999999
$ws.Range("L77").Value = 1473.3333
$ws.Range("M77").Value = -5188.666999999999
$ws.Range("N77").Value = -10209.3333
$ws.Range("H110").Value = 2864.5
$ws.Range("I110").Value = 2864.5
$ws.Range("K110").Value = 2864.5
$ws.Range("M110").Value = -819.5
$ws.Range("H132").Value = 2146.5454
$ws.Range("I132").Value = 2295.7778
$ws.Range("K132").Value = 6887.3334
$ws.Range("M132").Value = -4357.3334
$ws.Range("H136").Value = 2280.2942
$ws.Range("I136").Value = 2277.7334
$ws.Range("K136").Value = 6833.2002
$ws.Range("M136").Value = -4283.2002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1713.6666
$ws.Range("I20").Value = 738.7143
$ws.Range("K20").Value = 738.7143
$ws.Range("M20").Value = -491.7143
$ws.Range("H99").Value = 1985.5714
$ws.Range("I99").Value = 1985.5714
$ws.Range("K99").Value = 1985.5714
$ws.Range("M99").Value = -487.5714
$ws.Range("H105").Value = 1616.8
$ws.Range("I105").Value = 1616.8
$ws.Range("K105").Value = 1616.8
$ws.Range("M105").Value = 130.2
$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 55.6
$ws.Range("I7").Value = 55
$ws.Range("J7").Value = 56.5
$ws.Range("K7").Value = 55
$ws.Range("L7").Value = 56.5
$ws.Range("M7").Value = 58
$ws.Range("N7").Value = -282.5
$ws.Range("H56").Value = 35000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H58").Value = 2332
$ws.Range("I58").Value = 1798.2
$ws.Range("J58").Value = 3666.5
$ws.Range("K58").Value = 1798.2
$ws.Range("L58").Value = 3666.5
$ws.Range("M58").Value = -1595.2
$ws.Range("N58").Value = -4072.5
$ws.Range("H122").Value = 7719.8887
$ws.Range("I122").Value = 7997.375
$ws.Range("K122").Value = 23992.125
$ws.Range("M122").Value = -21542.125
$ws.Range("H132").Value = 3916
$ws.Range("I132").Value = 4374.75
$ws.Range("K132").Value = 13124.25
$ws.Range("M132").Value = -10594.25
$ws.Range("H136").Value = 2332
$ws.Range("I136").Value = 1798.2
$ws.Range("J136").Value = 3666.5
$ws.Range("K136").Value = 5394.6
$ws.Range("L136").Value = 10999.5
$ws.Range("M136").Value = -2844.6
$ws.Range("N136").Value = -16099.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 19799.8
$ws.Range("J106").Value = 19799.8
$ws.Range("L106").Value = 59399.39999999999
$ws.Range("N106").Value = -61291.39999999999
$ws.Range("H113").Value = 1194.5714
$ws.Range("J113").Value = 1259.4445
$ws.Range("L113").Value = 3778.3335
$ws.Range("N113").Value = -8118.333500000001
$ws.Range("H128").Value = 550996.75
$ws.Range("I128").Value = 550996.75
$ws.Range("K128").Value = 1652990.25
$ws.Range("M128").Value = -1648010.25
$ws.Range("H132").Value = 4300
$ws.Range("I132").Value = 3833.3333
$ws.Range("K132").Value = 34499.9997
$ws.Range("M132").Value = -31969.9997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 27999.75
$ws.Range("J39").Value = 27999.75
$ws.Range("L39").Value = 27999.75
$ws.Range("N39").Value = -29063.75
$ws.Range("H132").Value = 10452.643
$ws.Range("I132").Value = 10452.643
$ws.Range("K132").Value = 31357.929
$ws.Range("M132").Value = -28827.929
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H40").Value = 3721
$ws.Range("I40").Value = 3016.889
$ws.Range("K40").Value = 3016.889
$ws.Range("M40").Value = -2880.889
$ws.Range("H122").Value = 5751.9
$ws.Range("I122").Value = 4631.3887
$ws.Range("J122").Value = 7432.6665
$ws.Range("K122").Value = 13894.1661
$ws.Range("L122").Value = 22297.9995
$ws.Range("M122").Value = -11444.1661
$ws.Range("N122").Value = -27197.9995
$ws.Range("H136").Value = 4798.5
$ws.Range("I136").Value = 4798.5
$ws.Range("K136").Value = 14395.5
$ws.Range("M136").Value = -11845.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 39500
$ws.Range("J82").Value = 39500
$ws.Range("L82").Value = 39500
$ws.Range("N82").Value = -40266
$ws.Range("H85").Value = 39500
$ws.Range("J85").Value = 39500
$ws.Range("L85").Value = 39500
$ws.Range("N85").Value = -42152
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180
$ws.Range("H122").Value = 1531.5
$ws.Range("I122").Value = 1461.7273
$ws.Range("K122").Value = 4385.1819
$ws.Range("M122").Value = -1935.1819
$ws.Range("H124").Value = 15000
$ws.Range("J124").Value = 15000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -24820
$ws.Range("H132").Value = 2520.2856
$ws.Range("I132").Value = 2520.2856
$ws.Range("K132").Value = 7560.8568
$ws.Range("M132").Value = -5030.8568
$ws.Range("H136").Value = 2558.3684
$ws.Range("I136").Value = 2695
$ws.Range("J136").Value = 1397
$ws.Range("K136").Value = 8085
$ws.Range("L136").Value = 4191
$ws.Range("M136").Value = -5535
$ws.Range("N136").Value = -9291
